# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N ("Late"), shifting the old N/O/P ("Late" / heading / "Outstanding")
# columns one place to the right (-> O/P/Q). This is exactly what Excel's
# "Insert Column" does when a whole column is selected and Insert is invoked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate() | Out-Null

# Insert a new blank column at N; everything from N onward (N, O, P) shifts
# right by one (N->O, O->P, P->Q), matching the diff exactly.
$ws.Columns("N").EntireColumn.Insert() | Out-Null

# Give the newly inserted column N the same width as column M (its left
# neighbour) instead of the sheet's default width - this is the width Excel
# assigns to a freshly inserted column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Leave the selection where the author ended up after making the edit.
$ws.Range("R7").Select() | Out-Null
